# Compressor lab - progress update
# Updates measured velocity-head (B column) values across the water
# inlet/outlet sheets (now scaled into the correct units), refreshes the
# venturi-flow formula on water_50_outlet, tidies the data-row formatting
# (drop the heavy bottom border so the numbers sit flush right without the
# thick rule), and leaves the selection where each sheet was last clicked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# preface
# ---------------------------------------------------------------------
$wsPreface = $wb.Worksheets.Item("preface")
$wsPreface.Range("B7").Select()

# ---------------------------------------------------------------------
# water_25_outlet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("water_25_outlet")

$ws.Range("B2").Value = 8.6312223929999998
$ws.Range("B3").Value = 8.4126749610000005
$ws.Range("B4").Value = 7.957593857
$ws.Range("B5").Value = 7.474857858
$ws.Range("B6").Value = 6.1031958839999998
$ws.Range("B7").Value = 5.4588643509999999
$ws.Range("B8").Value = 4.7275152040000004
$ws.Range("B9").Value = 4.315611197

$rng = $ws.Range("A2:D9")
$rng.HorizontalAlignment = -4152
$rng.WrapText = $true
$rng.Borders.LineStyle = 0
$ws.Rows("2:9").EntireRow.AutoFit()

$ws.Range("D14").Select()

# ---------------------------------------------------------------------
# water_50_inlet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("water_50_inlet")

$ws.Range("B2").Value = 27.83485585
$ws.Range("B3").Value = 26.950983279999999
$ws.Range("B4").Value = 25.965494410000002
$ws.Range("B5").Value = 21.92058621
$ws.Range("B6").Value = 19.682215320000001
$ws.Range("B7").Value = 16.03178406
$ws.Range("B8").Value = 13.647160879999999
$ws.Range("B9").Value = 12.80217169
$ws.Range("B10").Value = 9.65
$ws.Range("B11").Value = 6.1031958839999998
$ws.Range("B12").Value = 3.86
$ws.Range("B13").Value = 2.7294321749999999

$rng = $ws.Range("A2:D13")
$rng.HorizontalAlignment = -4152
$rng.WrapText = $true
$rng.Borders.LineStyle = 0
$ws.Rows("2:13").EntireRow.AutoFit()

$ws.Range("G13").Select()

# ---------------------------------------------------------------------
# water_25_inlet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("water_25_inlet")

$ws.Range("B2").Value = 8.4126749610000005
$ws.Range("B3").Value = 7.957593857
$ws.Range("B4").Value = 7.474857858
$ws.Range("B5").Value = 6.958713962
$ws.Range("B6").Value = 6.4010858449999999
$ws.Range("B7").Value = 5.79
$ws.Range("B8").Value = 5.4588643509999999
$ws.Range("B9").Value = 4.315611197
$ws.Range("B10").Value = 3.86
$ws.Range("B11").Value = 3.3428580590000001
$ws.Range("B12").Value = 2.7294321749999999
$ws.Range("B13").Value = 1.93

$rng = $ws.Range("A2:D13")
$rng.HorizontalAlignment = -4152
$rng.WrapText = $true
$rng.Borders.LineStyle = 0
$ws.Rows("2:13").EntireRow.AutoFit()

$ws.Range("G4").Select()

# ---------------------------------------------------------------------
# water_50_outlet (select this one last so it ends up the active tab,
# matching the workbook's unchanged activeTab="1")
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("water_50_outlet")

$ws.Range("B2").Value = 27.76786452
$ws.Range("B3").Value = 27.430453880000002
$ws.Range("B4").Value = 26.742864470000001
$ws.Range("B5").Value = 25.965494410000002
$ws.Range("B6").Value = 24.182830689999999
$ws.Range("B7").Value = 21.317546759999999
$ws.Range("B8").Value = 19.10602523
$ws.Range("B9").Value = 15.67939412
$ws.Range("B10").Value = 12.80217169
$ws.Range("B11").Value = 9.2559548399999994
$ws.Range("B12").Value = 6.4010858449999999
$ws.Range("B13").Value = 3.3428580590000001

$ws.Range("F2").Formula = "=B2/(998.2 * 0.0064 * 0.0252)"

$rng = $ws.Range("A2:D13")
$rng.HorizontalAlignment = -4152
$rng.WrapText = $true
$rng.Borders.LineStyle = 0
$ws.Rows("2:13").EntireRow.AutoFit()

$ws.Range("L19").Select()
